# This script restores the original row order of the transaction tables on
# the "PayPal" and "eBay" worksheets. The source data had been shuffled;
# here we permute the data rows (columns A:D) back into the correct order
# for each sheet, row by row, using the Excel COM object model.

$wb = $excel.ActiveWorkbook

# Mapping of destination row number -> source row number (data rows only,
# the header row 1 is untouched) for the "PayPal" sheet.
$paypalMap = @{2=4; 3=16; 4=17; 5=5; 6=3; 7=11; 8=10; 9=13; 10=6; 11=12; 12=7; 13=9; 14=14; 15=8; 16=18; 17=2; 18=15}

# Mapping of destination row number -> source row number for the "eBay" sheet.
$ebayMap = @{2=23; 3=14; 4=13; 5=20; 6=2; 7=21; 8=7; 9=22; 10=15; 11=26; 12=17; 13=9; 14=4; 15=18; 16=5; 17=10; 18=6; 19=19; 20=8; 21=12; 22=25; 23=3; 24=24; 25=11; 26=16}

function Reorder-Sheet {
    param(
        $Workbook,
        $SheetName,
        $RowMap
    )

    $ws = $Workbook.Worksheets.Item($SheetName)

    # Snapshot every source row's values (columns A-D) before overwriting
    # anything, since several destination rows pull from rows that are also
    # being overwritten elsewhere in the table.
    $snapshot = @{}
    foreach ($srcRow in $RowMap.Values) {
        if (-not $snapshot.ContainsKey($srcRow)) {
            $rowValues = @(
                $ws.Cells.Item($srcRow, 1).Value2,
                $ws.Cells.Item($srcRow, 2).Value2,
                $ws.Cells.Item($srcRow, 3).Value2,
                $ws.Cells.Item($srcRow, 4).Value2
            )
            $snapshot[$srcRow] = $rowValues
        }
    }

    foreach ($destRow in $RowMap.Keys) {
        $srcRow = $RowMap[$destRow]
        $vals = $snapshot[$srcRow]

        # Force the cells to stay plain text (matching the source workbook,
        # where every cell in these tables is an inline/text string) instead
        # of letting Excel auto-detect dates or numbers from strings such as
        # "12/05/2023" or "-107.95".
        $rowRange = $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, 4))
        $rowRange.NumberFormat = "@"

        $ws.Cells.Item($destRow, 1).Value2 = $vals[0]
        $ws.Cells.Item($destRow, 2).Value2 = $vals[1]
        $ws.Cells.Item($destRow, 3).Value2 = $vals[2]
        $ws.Cells.Item($destRow, 4).Value2 = $vals[3]
    }
}

Reorder-Sheet $wb "PayPal" $paypalMap
Reorder-Sheet $wb "eBay" $ebayMap
